$d = $word.ActiveDocument

# The document's "DejaVu Sans" East-Asian (w:eastAsia) font is being
# switched to "Tahoma" on the styles that carry an explicit override:
# the Normal style and the Heading style (Word's Font.NameFarEast
# property is what round-trips to <w:rFonts w:eastAsia="...">).
$d.Styles("Normal").Font.NameFarEast = "Tahoma"
$d.Styles("Heading").Font.NameFarEast = "Tahoma"

# The List, Caption and Index styles previously had no rFonts override
# at all (empty <w:rPr/>) and now gain an explicit complex-script
# (w:cs) font of "DejaVu Sans". Font.NameBi is the Word object-model
# property that corresponds to the w:cs attribute of w:rFonts.
$d.Styles("List").Font.NameBi = "DejaVu Sans"
$d.Styles("Caption").Font.NameBi = "DejaVu Sans"
$d.Styles("Index").Font.NameBi = "DejaVu Sans"
